$d = $word.ActiveDocument

# The original last paragraph contains the placeholder {{paragraph_emphasis}}.
# It is replaced (split into three runs + proofErr spell-check markers) and
# followed by a large block of new paragraphs (list-template demo content).
$target = $d.Paragraphs.Last
$r = $target.Range

$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>paragraph_emphasis</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>}}</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:highlight w:val="yellow"/><w:lang w:val="en-US"/></w:rPr><w:t>Document Title: List Template Demo</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">  </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">This document demonstrates list patching capabilities.  </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-MX"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t>Ex</w:t></w:r><w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t>ample</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-MX"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t>simple_numbered</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t xml:space="preserve">}}  </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-MX"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t xml:space="preserve">Prueba de Viñetas Simples:  </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-MX"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t>simple_bullets</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t xml:space="preserve">}}  </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-MX"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t xml:space="preserve">  </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-MX"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t xml:space="preserve">Prueba de Viñetas Anidadas:  </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-MX"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t>nested_bullets</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t xml:space="preserve">}}  </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-MX"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t xml:space="preserve">  </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-MX"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t xml:space="preserve">Prueba de Numeración Anidada:  </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-MX"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t>nested_numbered</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t xml:space="preserve">}}  </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-MX"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t xml:space="preserve">  </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-MX"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t xml:space="preserve">Prueba de Formato Complejo:  </w:t></w:r></w:p><w:p><w:r><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>complex_formatting</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}}</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:lang w:val="es-MX"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t xml:space="preserve">Prueba de Lista de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t>checkbox</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t>checklist</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t>)</w:t></w:r></w:p><w:p><w:r><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>checkbox_list</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}}</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>
'@

[void]$r.InsertXML($newXml)
